$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.460.57'
$ws.Range("E2").Value = '  +0.76%  '

# Row 3
$ws.Range("D3").Value = '3.692.50'
$ws.Range("E3").Value = '  +1.34%  '

# Row 4
$ws.Range("E4").Value = '  -0.71%  '

# Row 5
$ws.Range("D5").Value = '''686.72'
$ws.Range("E5").Value = '  +2.90%  '

# Row 6
$ws.Range("D6").Value = '''160.90'
$ws.Range("E6").Value = '  +1.61%  '

# Row 7
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.59%  '

# Row 8
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = '''0.495'
$ws.Range("E8").Value = '  +3.42%  '

# Row 9
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '''0.147'
$ws.Range("E9").Value = '  +1.26%  '

# Row 10
$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").Value = '''7.14'
$ws.Range("E10").Value = '  +3.48%  '

# Row 11
$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").Value = '''0.436'
$ws.Range("E11").Value = '  -0.70%  '

# Row 12
$ws.Range("B12").Value = 'ShibaInu'
$ws.Range("C12").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D12").Value = '''0.0000234'
$ws.Range("E12").Value = '  +5.46%  '

# Row 13
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '4.315.40'
$ws.Range("E13").Value = '  +0.92%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '''32.64'
$ws.Range("E14").Value = '  +1.39%  '

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.695.72'
$ws.Range("E15").Value = '  +0.01%  '

# Row 16
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '69.473.12'
$ws.Range("E16").Value = '  +1.15%  '

# Row 17
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '''0.115'
$ws.Range("E17").Value = '  +1.65%  '

# Row 18
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '''15.91'
$ws.Range("E18").Value = '  +0.82%  '

# Row 19
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '''6.48'
$ws.Range("E19").Value = '  +2.39%  '

# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''472.21'
$ws.Range("E20").Value = '  -0.48%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''10.07'
$ws.Range("E21").Value = '  +6.32%  '

# Row 22
$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D22").Value = '''0.652'
$ws.Range("E22").Value = '  +1.86%  '

# Row 23
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '''79.75'
$ws.Range("E23").Value = '  -0.69%  '

# Row 24
$ws.Range("B24").Value = 'WrappedeETH'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D24").Value = '3.835.29'
$ws.Range("E24").Value = '  +0.14%  '

# Row 25
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  +0.53%  '

# Row 26
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").Value = '''0.0000125'
$ws.Range("E26").Value = '  +5.11%  '

# Row 27
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '''11.06'
$ws.Range("E27").Value = '  -1.51%  '

# Row 28
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''9.26'
$ws.Range("E28").Value = '  +4.16%  '

# Row 29
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''2.72'
$ws.Range("E29").Value = '  +0.66%  '

# Row 30
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").Value = '''1.73'
$ws.Range("E30").Value = '  -0.22%  '

# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '''2.02'
$ws.Range("E31").Value = '  +1.53%  '

# Row 32
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = '''6.61'
$ws.Range("E32").Value = '  +2.64%  '

# Row 33
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''27.01'
$ws.Range("E33").Value = '  +3.93%  '

# Row 34
$ws.Range("D34").Value = '''0.998'
$ws.Range("E34").Value = '  -0.38%  '

# Row 35
$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").Value = '3.664.95'
$ws.Range("E35").Value = '  +1.15%  '

# Row 36
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '''0.161'
$ws.Range("E36").Value = '  +2.77%  '

# Row 37
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = '''8.21'
$ws.Range("E37").Value = '  +0.09%  '

# Row 38
$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").Value = '''6.18'
$ws.Range("E38").Value = '  +5.39%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '''2.24'
$ws.Range("E39").Value = '  +4.87%  '

# Row 40
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '  +0.01%  '

# Row 41
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '''0.0907'
$ws.Range("E41").Value = '  +2.94%  '

# Row 42
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  -0.52%  '

# Row 43
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = '''0.944'
$ws.Range("E43").Value = '  +1.56%  '

# Row 44
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").Value = '''165.48'
$ws.Range("E44").Value = '  +5.58%  '

# Row 45
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '''47.42'
$ws.Range("E45").Value = '  -1.25%  '

# Row 46
$ws.Range("B46").Value = 'SuiNetwork'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D46").Value = '''1.14'
$ws.Range("E46").Value = '  +11.45%  '

# Row 47
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = '''2.75'
$ws.Range("E47").Value = '  +0.67%  '

# Row 48
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Value = '''0.000283'
$ws.Range("E48").Value = '  +7.36%  '

# Row 49
$ws.Range("D49").Value = '''1.32'
$ws.Range("E49").Value = '  +3.75%  '

# Row 50
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''28.11'
$ws.Range("E50").Value = '  +3.58%  '

# Row 51
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '''7.85'
$ws.Range("E51").Value = '  +1.09%  '
